# LoginTest.xlsx edit: remove the "blank_password_error_xpath" CHECKTEXT row
# from the "BlankInput" sheet (4th sheet). This drops the shared strings
# "Can't be blank." and "blank_password_error_xpath" (no longer referenced
# anywhere in the workbook) and shifts the following row (QUIT/NA/NA/NA) up
# to take its place.

$wb = $excel.ActiveWorkbook

# "BlankInput" is the 4th sheet in the workbook.
$ws = $wb.Worksheets.Item("BlankInput")

# Row 6 currently holds: CHECKTEXT | blank_password_error_xpath | XPATH | Can't be blank.
# Deleting it shifts the old row 7 (QUIT | NA | NA | NA) up into row 6.
$ws.Rows.Item(6).Delete() | Out-Null

# Update the sheet's active selection to reflect the new last used cell.
$ws.Range("A9").Select() | Out-Null
